# Criação dos cadastros de status e local de tratamento
# Adds new GRANT rows (local/status de tratamento + risco do paciente) to the
# "grants por usuario" sheet, replacing the stray formatted-but-empty row 32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# The sheet currently ends with an empty, oddly-styled row 32 (A32 only,
# 14pt Times New Roman). Remove it first so the new data rows can occupy
# rows 30-35 contiguously after row 29.
$ws.Rows.Item(32).Delete()

# Prime rows 30:35 with the same formatting as row 29 (font/alignment on
# column B, plain text elsewhere) so new cells don't invent new styles.
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New grant statements: tratamento.tb_c_local_trtmto / sq_local_trtmto,
# tb_c_status_trtmto / sq_status_trtmto, tb_c_risco_pcnt / sq_risco_pcnt.
$grants = @(
    @("GRANT SELECT, UPDATE, INSERT, DELETE ON tratamento.tb_c_local_trtmto TO ", "evaldo", ";"),
    @("GRANT SELECT ON SEQUENCE tratamento.sq_local_trtmto TO", "evaldo", ";"),
    @("GRANT SELECT, UPDATE, INSERT, DELETE ON tratamento.tb_c_status_trtmto TO ", "evaldo", ";"),
    @("GRANT SELECT ON SEQUENCE tratamento.sq_status_trtmto TO", "evaldo", ";"),
    @("GRANT SELECT, UPDATE, INSERT, DELETE ON tratamento.tb_c_risco_pcnt TO ", "evaldo", ";"),
    @("GRANT SELECT ON SEQUENCE tratamento.sq_risco_pcnt TO", "evaldo", ";")
)

$r = 30
foreach ($row in $grants) {
    $ws.Range("A$r").Value2 = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $ws.Range("D$r").Formula = "=A$r&"" ""&B$r&"" ""&C$r"
    $r++
}

# Column D previously was sized to fit the longest existing string; after
# the edit the widest entry is shorter, so the best-fit width shrinks.
$ws.Columns.Item(4).ColumnWidth = 99.3

# Clear the stale D2:D29 selection left over from before the rows were
# appended.
$ws.Range("A1").Select()
